$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("metro_budget")

$ws.Range("F2").Formula = '=IFERROR(RANK(E2,$E$2:$E$52,1),"-")'
$ws.Range("K2").Formula = '=IFERROR(RANK(J2,$J$2:$J$52,1),"-")'
$ws.Range("P2").Formula = '=IFERROR(RANK(O2,$O$2:$O$52,1),"-")'

$ws.Range("F3:F52").Formula = '=IFERROR(RANK(E3,$E$2:$E$52,1),"-")'
$ws.Range("K3:K52").Formula = '=IFERROR(RANK(J3,$J$2:$J$52,1),"-")'
$ws.Range("P3:P52").Formula = '=IFERROR(RANK(O3,$O$2:$O$52,1),"-")'

$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Application.ActiveWindow.ScrollColumn = 6
$ws.Range("P7").Select()
